$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "Africana A/V"
$ws.Range("E8").Value = "Africana Library > A/V"
$ws.Range("D12").Select() | Out-Null
